# Practical Exercise 5 dan 6
# Adds firstName / lastName / postalCode columns (D:F) with sample data,
# matching the sharedStrings / sheet1 diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1)
$ws.Range("D1").Value = "firstName"
$ws.Range("E1").Value = "lastName"
$ws.Range("F1").Value = "postalCode"

# New data row (row 2)
$ws.Range("D2").Value = "Cintya"
$ws.Range("E2").Value = "Nainggolan"
$ws.Range("F2").Value = 12020

# Column widths: D -> 13, E:F -> 11 (ColumnWidth property already excludes
# the fixed cell-padding Excel adds internally, so back it off here so the
# stored <col> width lands exactly on 13 / 11).
$ws.Columns.Item(4).ColumnWidth = 12.1666666666667
$ws.Columns.Item(5).ColumnWidth = 10.1666666666667
$ws.Columns.Item(6).ColumnWidth = 10.1666666666667

# Move the active selection to F2, matching the saved sheet view.
[void]$ws.Range("F2").Select()
